$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Target table (row, A=countdown number, B=From, C=To, D=Drone) — columns
# B/C/D hold digit-strings that must stay TEXT (shared strings), matching the
# original sheet's convention, not be auto-coerced to numbers by Excel.
# ---------------------------------------------------------------------------
$data = @(
  @(2,  9, "1", "7", "1"),
  @(3,  8, "7", "1", "1"),
  @(4,  7, "2", "9", "12"),
  @(5,  6, "9", "2", "12"),
  @(6,  5, "1", "6", "13"),
  @(7,  4, "3", "5", "13"),
  @(8,  3, "4", "8", "13"),
  @(9,  2, "5", "4", "13"),
  @(10, 1, "6", "3", "13"),
  @(11, 0, "8", "1", "13")
)

# A scratch cell used purely to stamp a TEXT number-format onto a value so
# that, once copied, the digits paste into the destination as a text string
# instead of being re-parsed as a number (Excel's normal autoconvert
# behaviour for bare digit strings assigned via .Value).
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

foreach ($row in $data) {
  $r = $row[0]

  # Column A: plain numeric value.
  $ws.Cells.Item($r, 1).Value = $row[1]

  # Columns B, C, D: force text/shared-string storage via the scratch cell.
  for ($col = 2; $col -le 4; $col++) {
    $scratch.Value = $row[$col]
    $scratch.Copy()
    $ws.Cells.Item($r, $col).PasteSpecial(-4163)
  }
}

# Column A on the newly-added rows (8-11) needs the same bold/bordered,
# centered style already used by A2:A7. Copy formats only so no brand new
# style entry is introduced.
$ws.Range("A7").Copy()
$ws.Range("A8:A11").PasteSpecial(-4122)

# Tidy up the scratch cell so it doesn't leak into the saved sheet.
$scratch.Clear()

$excel.CutCopyMode = 0
